$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results (A2:A25 hold case indices 0..23).
# Columns D,H,I,K,O are untouched (still 0); B,C,E,F,G,J,L,M,N get new values.
$newValues = @{
    2 = @{ "B"=26.93347275024401; "C"=10.79180077743248; "E"=9.530651036881768; "F"=51.02769897572396; "G"=3.767228628479839; "J"=10.66347546743639; "L"=10.80256174471938; "M"=20.62101746361374; "N"=22.88875440455603 }
    3 = @{ "B"=26.62515783084871; "C"=10.30293347031096; "E"=9.460399037181066; "F"=50.90712187817036; "G"=3.771850442160405; "J"=10.6893809605731; "L"=10.8203071916697; "M"=20.58735820633351; "N"=22.93784059433014 }
    4 = @{ "B"=26.44221883679043; "C"=9.994378523491442; "E"=9.416199038851149; "F"=50.84441049598794; "G"=3.774833378818084; "J"=10.70612869831473; "L"=10.83239420775929; "M"=20.57115348027778; "N"=22.96992686155001 }
    5 = @{ "B"=26.36935393892921; "C"=9.866736563555282; "E"=9.397921180184756; "F"=50.82171097140979; "G"=3.776085592631694; "J"=10.71316562951734; "L"=10.83761989845751; "M"=20.56567587304203; "N"=22.9834915746032 }
    6 = @{ "B"=26.35735885238991; "C"=9.84543330490251; "E"=9.394870142403978; "F"=50.8181144146873; "G"=3.776295739328903; "J"=10.7143469297319; "L"=10.83850576477886; "M"=20.56483442165665; "N"=22.98577352809889 }
    7 = @{ "B"=26.44122922899375; "C"=9.99266451391061; "E"=9.415953611924294; "F"=50.84409278966512; "G"=3.774850118044952; "J"=10.70622274146708; "L"=10.83246346724138; "M"=20.57107504389581; "N"=22.97010781925922 }
    8 = @{ "B"=26.82589511325336; "C"=10.62509581743439; "E"=9.506650253976288; "F"=50.98377971199214; "G"=3.768792199911168; "J"=10.67223326543246; "L"=10.80843348701824; "M"=20.60848821627461; "N"=22.90527499223661 }
    9 = @{ "B"=27.62693182381649; "C"=11.79098467985171; "E"=9.675981763594885; "F"=51.34706404542688; "G"=3.758057283711274; "J"=10.61223760218894; "L"=10.77073694219101; "M"=20.7170536296182; "N"=22.79360237746328 }
    10 = @{ "B"=28.23863797241576; "C"=12.5935457802159; "E"=9.795101402278695; "F"=51.66762748690701; "G"=3.750858548742052; "J"=10.57218852038644; "L"=10.74875304901649; "M"=20.81794478504025; "N"=22.72100407989075 }
    11 = @{ "B"=28.52079796867522; "C"=12.9455003681062; "E"=9.848124927152341; "F"=51.82488190136468; "G"=3.74773103992043; "J"=10.55483814706706; "L"=10.7399847408086; "M"=20.86834234086541; "N"=22.69003356602368 }
    12 = @{ "B"=28.62810700247421; "C"=13.07679296208607; "E"=9.868034344039726; "F"=51.88604961903773; "G"=3.746567746492909; "J"=10.54839242431503; "L"=10.73684097168885; "M"=20.88806501059467; "N"=22.67860184372635 }
    13 = @{ "B"=28.60497713176963; "C"=13.04860640500234; "E"=9.863754061906731; "F"=51.8728044733043; "G"=3.746817349709152; "J"=10.54977509492684; "L"=10.73751019402155; "M"=20.88378915030578; "N"=22.68105068792352 }
    14 = @{ "B"=28.52961766656733; "C"=12.95634216731745; "E"=9.849766293890017; "F"=51.82988191084718; "G"=3.747634914504901; "J"=10.55430536065355; "L"=10.73972256520535; "M"=20.86995220131635; "N"=22.68908713239769 }
    15 = @{ "B"=28.48351496455659; "C"=12.89956656275704; "E"=9.841176257878471; "F"=51.8038006577738; "G"=3.748138430477678; "J"=10.55709648082945; "L"=10.74110068681137; "M"=20.86155949925636; "N"=22.6940482728378 }
    16 = @{ "B"=28.22026803571573; "C"=12.57027183564453; "E"=9.791612586091759; "F"=51.65757847338248; "G"=3.751065889037473; "J"=10.5733398471923; "L"=10.74935082400083; "M"=20.81474101022304; "N"=22.72306947956231 }
    17 = @{ "B"=28.0597011400902; "C"=12.36482261072143; "E"=9.760907183137066; "F"=51.57078637283094; "G"=3.752899396049463; "J"=10.58352672725649; "L"=10.75472720716787; "M"=20.7871661408311; "N"=22.74139979626398 }
    18 = @{ "B"=27.96771950425046; "C"=12.24542205593373; "E"=9.74313674630652; "F"=51.52194318864887; "G"=3.753967847626001; "J"=10.58946768650555; "L"=10.75793558476426; "M"=20.77173012319677; "N"=22.7521362772851 }
    19 = @{ "B"=27.93664303660755; "C"=12.20478674315353; "E"=9.737101235533226; "F"=51.50559146466514; "G"=3.754331993048679; "J"=10.59149324223611; "L"=10.75904183167842; "M"=20.76657689415454; "N"=22.75580465674192 }
    20 = @{ "B"=28.07675595460184; "C"=12.38682124340328; "E"=9.764187173054101; "F"=51.5799142094934; "G"=3.752702781987336; "J"=10.58243385978595; "L"=10.75414287825698; "M"=20.7900576809278; "N"=22.73942848527618 }
    21 = @{ "B"=28.55174083361924; "C"=12.98349697292226; "E"=9.853879454571334; "F"=51.84244557968965; "G"=3.74739420641677; "J"=10.55297133600337; "L"=10.73906795015977; "M"=20.87399920464363; "N"=22.68671859129842 }
    22 = @{ "B"=28.86481389651867; "C"=13.36185492640488; "E"=9.911510023418794; "F"=52.02344798117482; "G"=3.744047235366625; "J"=10.53444132369633; "L"=10.73024468020184; "M"=20.93257500575739; "N"=22.65399608286926 }
    23 = @{ "B"=28.69751139848176; "C"=13.16100785990781; "E"=9.880842596680592; "F"=51.92599013285001; "G"=3.745822416959715; "J"=10.54426487484107; "L"=10.73485986705052; "M"=20.90097528673112; "N"=22.67130250453974 }
    24 = @{ "B"=28.06904443632781; "C"=12.37687965677591; "E"=9.76270465542239; "F"=51.57578422739593; "G"=3.752791626515532; "J"=10.58292768229817; "L"=10.7544066875648; "M"=20.78874911566338; "N"=22.74031909786826 }
    25 = @{ "B"=27.40576869450791; "C"=11.48452094539608; "E"=9.631095175337485; "F"=51.23929303834558; "G"=3.760839814373407; "J"=10.62775830454165; "L"=10.7799292652223; "M"=20.68394610537724; "N"=22.82215524916465 }
}

foreach ($row in $newValues.Keys) {
    $rowData = $newValues[$row]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
